# "jogos, batedores e noticias"
# Update the "M. United" batedores entries: the shared string
# "CR7 - Bruno Fernandes" is replaced by two distinct strings,
# "Bruno Fernandes" (Pênaltis column) and " Bruno Fernandes" (Faltas column,
# with a leading space), since Cristiano Ronaldo left Manchester United.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("C16").Value = "Bruno Fernandes"
$ws.Range("D16").Value = " Bruno Fernandes"

# Restore the selection/active cell to D17 (previously D21) as left by the author.
$ws.Range("D17").Select()
